# Apply cryptos list update (Mon Sep 30 14:31:12 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.690.07"
$ws.Range("E2").Value = "  -3.11%  "
$ws.Range("D3").Value = "2.609.34"
$ws.Range("E3").Value = "  -1.68%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.66"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.21"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.64%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -1.60%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.119"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -5.77%  "
$ws.Range("E10").Value = "  -0.17%  "
$ws.Range("E11").Value = "  -3.89%  "
$ws.Range("E12").Value = "  -0.39%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.10"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.74%  "
$ws.Range("D14").Value = "3.082.64"
$ws.Range("E14").Value = "  -1.58%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000182"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").Value = "63.476.64"
$ws.Range("E16").Value = "  -3.19%  "
$ws.Range("D17").Value = "2.607.14"
$ws.Range("E17").Value = "  +0.09%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -4.08%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.60"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("E20").Value = "  -4.16%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "342.05"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.19%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.08"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.74"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.04%  "
$ws.Range("E25").Value = "  -1.77%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "594.06"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.91%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.17"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -4.35%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.57"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.46%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.03%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.160"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -2.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.85"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -3.10%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.72"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.53"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  -2.36%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.408"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -2.71%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "19.72"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.998"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "154.95"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.86"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.01%  "
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("E42").Value = "  -2.63%  "
$ws.Range("E43").Value = "  +6.26%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "155.69"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.62%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "3.92"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -3.83%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.88"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0591"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.50%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.101"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +0.95%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.626"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.88%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0247"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -3.36%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "18.89"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -4.21%  "
